$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings are stored as text (matches source inlineStr cells)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.258.45"
$ws.Range("E2").Value = "  -0.82%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.082.36"
$ws.Range("E3").Value = "  +3.30%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9989"
$ws.Range("E4").Value = "  -0.31%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "327.51"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9970"
$ws.Range("E6").Value = "  -0.35%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5188"
$ws.Range("E7").Value = "  +1.45%  "
$ws.Range("E8").Value = "  +3.81%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08817"
$ws.Range("E9").Value = "  +0.37%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "46.10"
$ws.Range("E10").Value = "  +6.70%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.158"
$ws.Range("E11").Value = "  +1.85%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.35"
$ws.Range("E12").Value = "  -1.29%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.074.59"
$ws.Range("E13").Value = "  +3.28%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.666"
$ws.Range("E14").Value = "  +0.70%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.669"
$ws.Range("E15").Value = "  +2.19%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "94.99"
$ws.Range("E16").Value = "  +0.48%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9985"
$ws.Range("E17").Value = "  -0.48%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001119"
$ws.Range("E18").Value = "  +0.13%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06604"
$ws.Range("E19").Value = "  +1.02%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.77"
$ws.Range("E20").Value = "  -0.91%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9971"
$ws.Range("E21").Value = "  -0.32%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.290"
$ws.Range("E22").Value = "  +0.94%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "30.295.34"
$ws.Range("E23").Value = "  -0.93%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.24"
$ws.Range("E24").Value = "  +3.09%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.284"
$ws.Range("E25").Value = "  +2.49%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.313.05"
$ws.Range("E26").Value = "  +2.87%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.35"
$ws.Range("E27").Value = "  -0.24%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.583"
$ws.Range("E28").Value = "  +5.98%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "161.55"
$ws.Range("E29").Value = "  -1.04%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "130.76"
$ws.Range("E30").Value = "  -0.61%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.185"
$ws.Range("E31").Value = "  +3.92%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1064"
$ws.Range("E32").Value = "  +1.00%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.650"
$ws.Range("E33").Value = "  +21.66%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.150"
$ws.Range("E34").Value = "  +0.85%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.817"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02574"
$ws.Range("E36").Value = "  +1.89%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "9.762"
$ws.Range("E37").Value = "  +7.31%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "12.68"
$ws.Range("E38").Value = "  +2.84%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2243"
$ws.Range("E41").Value = "  +1.99%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6780"
$ws.Range("E42").Value = "  +1.41%  "
$ws.Range("E43").Value = "  +1.14%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9967"
$ws.Range("E44").Value = "  -0.30%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.02"
$ws.Range("E45").Value = "  +2.82%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6325"
$ws.Range("E46").Value = "  +2.09%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.198"
$ws.Range("E47").Value = "  -0.17%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.600"
$ws.Range("E48").Value = "  -1.84%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.235"
$ws.Range("E49").Value = "  -2.81%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.191"
$ws.Range("E50").Value = "  +7.50%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "81.34"
$ws.Range("E51").Value = "  -0.03%  "

# Rows 39/40: Hedera and InternetComputer(DFINITY) swap positions, with refreshed price data
$ws.Range("D39").NumberFormat = "@"
$ws.Range("B39").Value = "InternetComputer(DFINITY)"
$ws.Range("C39").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D39").Value = "5.419"
$ws.Range("E39").Value = "  -0.80%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").Value = "0.06620"
$ws.Range("E40").Value = "  -0.77%  "
